$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F11").Value = 16
$ws.Range("G11").Value = 2329.44
$ws.Range("B14").Value = 19668.37
$ws.Range("F26").Value = 116
$ws.Range("G26").Value = 5942.68
$ws.Range("F33").Value = 13
$ws.Range("G33").Value = 998.92
$ws.Range("B36").Value = 13602.32
$ws.Range("F41").Value = 104
$ws.Range("G41").Value = 7284.16
$ws.Range("F51").Value = 11
$ws.Range("G51").Value = 278.19
$ws.Range("B71").Value = 78569.66
$ws.Range("F100").Value = 42
$ws.Range("G100").Value = 7225.68
$ws.Range("F108").Value = 56
$ws.Range("G108").Value = 5331.76
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("F116").Value = 91
$ws.Range("G116").Value = 12280.45
$ws.Range("F141").Value = 637
$ws.Range("G141").Value = 12396.02
$ws.Range("B143").Value = 355969.28
$ws.Range("F173").Value = 8
$ws.Range("G173").Value = 195.44
$ws.Range("B176").Value = 24189.58
$ws.Range("F193").Value = 107
$ws.Range("G193").Value = 5294.36
$ws.Range("F194").Value = 42
$ws.Range("G194").Value = 2078.16
$ws.Range("F199").Value = 308
$ws.Range("G199").Value = 6061.44
$ws.Range("F201").Value = 170
$ws.Range("G201").Value = 7568.4
$ws.Range("B205").Value = 35670.26
$ws.Range("B213").Value = 57756
$ws.Range("B214").Value = 53925
$ws.Range("F283").Value = 133
$ws.Range("G283").Value = 11282.39
$ws.Range("B288").Value = 24893.64
$ws.Range("F356").Value = 61
$ws.Range("G356").Value = 2296.65
$ws.Range("B371").Value = 156284.82
$ws.Range("F378").Value = 3
$ws.Range("G378").Value = 575.7
$ws.Range("F383").Value = 6
$ws.Range("G383").Value = 442.2
$ws.Range("F410").Value = 282
$ws.Range("G410").Value = 31572.72
$ws.Range("F419").Value = 63
$ws.Range("G419").Value = 6557.67
$ws.Range("F423").Value = 273
$ws.Range("G423").Value = 27600.3
$ws.Range("F430").Value = 127
$ws.Range("G430").Value = 15066.01
$ws.Range("F431").Value = 271
$ws.Range("G431").Value = 16024.23
$ws.Range("F436").Value = 786
$ws.Range("G436").Value = 17048.34
$ws.Range("F437").Value = 389
$ws.Range("G437").Value = 2337.89
$ws.Range("F438").Value = 33
$ws.Range("G438").Value = 2742.3
$ws.Range("F440").Value = 27
$ws.Range("G440").Value = 11418.3
$ws.Range("B456").Value = 661908.96
$ws.Range("F481").Value = 124
$ws.Range("G481").Value = 4476.4
$ws.Range("B489").Value = 7851.09
$ws.Range("F516").Value = 87
$ws.Range("G516").Value = 5972.55
$ws.Range("F525").Value = 250
$ws.Range("G525").Value = 24150
$ws.Range("B531").Value = 76694.88
$ws.Range("F546").Value = 43
$ws.Range("G546").Value = 12625.66
$ws.Range("B560").Value = 83127.09
$ws.Range("F587").Value = 1175
$ws.Range("G587").Value = 7731.5
$ws.Range("B594").Value = 124452.42
$ws.Range("F648").Value = 0
$ws.Range("G648").Value = 0
$ws.Range("B649").Value = 29698.25
$ws.Range("F662").Value = 564
$ws.Range("G662").Value = 3778.8
$ws.Range("B666").Value = 50381.89
$ws.Range("F710").Value = 32
$ws.Range("G710").Value = 114.88
$ws.Range("F716").Value = 61
$ws.Range("G716").Value = 2436.34
$ws.Range("B728").Value = 15214.66
$ws.Range("F731").Value = 46
$ws.Range("G731").Value = 13150.02
$ws.Range("F736").Value = 46
$ws.Range("G736").Value = 2310.58
$ws.Range("F746").Value = 91
$ws.Range("G746").Value = 6693.05
$ws.Range("B750").Value = 96348.17
$ws.Range("F752").Value = 137
$ws.Range("G752").Value = 17885.35
$ws.Range("F755").Value = 136
$ws.Range("G755").Value = 12913.2
$ws.Range("F758").Value = 153
$ws.Range("G758").Value = 4161.6
$ws.Range("B759").Value = 83913.42999999999
$ws.Range("F791").Value = 214
$ws.Range("G791").Value = 9240.52
$ws.Range("B793").Value = 56197.8
$ws.Range("F822").Value = 59
$ws.Range("G822").Value = 5413.25
$ws.Range("F827").Value = 120
$ws.Range("G827").Value = 9645.6
$ws.Range("B828").Value = 31449.95
$ws.Range("F835").Value = 24
$ws.Range("G835").Value = 14033.28
$ws.Range("F836").Value = 7
$ws.Range("G836").Value = 3744.16
$ws.Range("B840").Value = 45914.82
$ws.Range("F859").Value = 13
$ws.Range("G859").Value = 1332.63
$ws.Range("F860").Value = 5
$ws.Range("G860").Value = 8026.5
$ws.Range("F867").Value = 16
$ws.Range("G867").Value = 9972.48
$ws.Range("B869").Value = 68707.95
$ws.Range("F883").Value = 51
$ws.Range("G883").Value = 2740.23
$ws.Range("F886").Value = 96
$ws.Range("G886").Value = 10639.68
$ws.Range("F887").Value = 174
$ws.Range("G887").Value = 25056
$ws.Range("B890").Value = 138077.24
$ws.Range("F909").Value = 324
$ws.Range("G909").Value = 9758.879999999999
$ws.Range("B922").Value = 78379.69
$ws.Range("F927").Value = 95
$ws.Range("G927").Value = 622.25
$ws.Range("F928").Value = 259
$ws.Range("G928").Value = 9686.6
$ws.Range("F930").Value = 239
$ws.Range("G930").Value = 8938.6
$ws.Range("B931").Value = 33267.37
$ws.Range("F982").Value = 0
$ws.Range("G982").Value = 0
$ws.Range("F983").Value = 7
$ws.Range("G983").Value = 1678.53
$ws.Range("F984").Value = 8
$ws.Range("G984").Value = 1833.68
$ws.Range("B986").Value = 12367.08
$ws.Range("F991").Value = 112
$ws.Range("G991").Value = 16560.32
$ws.Range("B994").Value = 564363.1800000001
$ws.Range("F996").Value = 2
$ws.Range("G996").Value = 625.5
$ws.Range("B997").Value = 625.5
$ws.Range("B1001").Value = 4600512.18
$ws.Range("B1002").Value = 4600512.18
